$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot).
# Columns D (Price) and E (Volume 1h) are forced to Text format before
# assignment so numeric-looking strings (e.g. "0.9983", "29.788.76")
# are preserved exactly as scraped rather than being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.788.76'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +8.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.953.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4788'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4144'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +8.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.38'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08310'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.053'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +8.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.86'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.30%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.199'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.45%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.925.25'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.464'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '93.26'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9989'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06703'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.13'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9976'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.749.99'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.655'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.268'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.179.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +6.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.77'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.26'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.210'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.640'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.86'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.034'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +10.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09660'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.64%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +11.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.681'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.518'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.32%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06272'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.68%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02324'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.753'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.209'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6141'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.80'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1918'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9971'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.278'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.62'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5754'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.325'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +27.44%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07315'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.36'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.89%  '
